# Update gh-pages output (想去人数 / column F counts) across all sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 145
$ws.Cells.Item(7, 6).Value = 13312
$ws.Cells.Item(8, 6).Value = 347
$ws.Cells.Item(9, 6).Value = 5364
$ws.Cells.Item(11, 6).Value = 7221
$ws.Cells.Item(14, 6).Value = 3708
$ws.Cells.Item(18, 6).Value = 187
$ws.Cells.Item(22, 6).Value = 89
$ws.Cells.Item(24, 6).Value = 3799
$ws.Cells.Item(25, 6).Value = 119
$ws.Cells.Item(26, 6).Value = 5055
$ws.Cells.Item(28, 6).Value = 2020
$ws.Cells.Item(29, 6).Value = 119
$ws.Cells.Item(30, 6).Value = 310
$ws.Cells.Item(31, 6).Value = 7436
$ws.Cells.Item(32, 6).Value = 27
$ws.Cells.Item(33, 6).Value = 170
$ws.Cells.Item(34, 6).Value = 2114
$ws.Cells.Item(36, 6).Value = 138
$ws.Cells.Item(37, 6).Value = 1144
$ws.Cells.Item(40, 6).Value = 244
$ws.Cells.Item(41, 6).Value = 236
$ws.Cells.Item(43, 6).Value = 1166
$ws.Cells.Item(44, 6).Value = 1165
$ws.Cells.Item(45, 6).Value = 19
$ws.Cells.Item(47, 6).Value = 1279
$ws.Cells.Item(48, 6).Value = 1952
$ws.Cells.Item(49, 6).Value = 101
$ws.Cells.Item(50, 6).Value = 187

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 12
$ws.Cells.Item(5, 6).Value = 141
$ws.Cells.Item(10, 6).Value = 932
$ws.Cells.Item(12, 6).Value = 110

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 524
$ws.Cells.Item(3, 6).Value = 705
$ws.Cells.Item(4, 6).Value = 57

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 145
$ws.Cells.Item(5, 6).Value = 524
$ws.Cells.Item(6, 6).Value = 705
$ws.Cells.Item(7, 6).Value = 57
$ws.Cells.Item(9, 6).Value = 13311
$ws.Cells.Item(10, 6).Value = 347
$ws.Cells.Item(11, 6).Value = 5363
$ws.Cells.Item(12, 6).Value = 3708
$ws.Cells.Item(16, 6).Value = 187
$ws.Cells.Item(19, 6).Value = 89
$ws.Cells.Item(20, 6).Value = 12
$ws.Cells.Item(21, 6).Value = 141
$ws.Cells.Item(22, 6).Value = 3799
$ws.Cells.Item(24, 6).Value = 119
$ws.Cells.Item(25, 6).Value = 5055
$ws.Cells.Item(27, 6).Value = 2020
$ws.Cells.Item(28, 6).Value = 119
$ws.Cells.Item(29, 6).Value = 311
$ws.Cells.Item(30, 6).Value = 7436
$ws.Cells.Item(31, 6).Value = 27
$ws.Cells.Item(32, 6).Value = 170
$ws.Cells.Item(33, 6).Value = 2114
$ws.Cells.Item(35, 6).Value = 138
$ws.Cells.Item(36, 6).Value = 1144
$ws.Cells.Item(38, 6).Value = 244
$ws.Cells.Item(39, 6).Value = 236
$ws.Cells.Item(40, 6).Value = 1166
$ws.Cells.Item(41, 6).Value = 1165
$ws.Cells.Item(42, 6).Value = 19
$ws.Cells.Item(45, 6).Value = 1279
$ws.Cells.Item(47, 6).Value = 1952
$ws.Cells.Item(48, 6).Value = 101
$ws.Cells.Item(50, 6).Value = 187
